# Weekly update: a new "Poroto granado" price record (Región de O'Higgins,
# 2022-01-24) is inserted into the dataset right after row 359. This shifts
# every subsequent data row (old 360..381) down by one (new 361..382),
# preserving all of their existing values, and the new row 360 receives the
# fresh record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 360, pushing rows 360:381 down to 361:382.
$ws.Rows.Item(360).Insert()

$newRow = 360
$ws.Cells.Item($newRow, 1).Value = 6
$ws.Cells.Item($newRow, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($newRow, 3).Value = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value = 44585
$ws.Cells.Item($newRow, 5).Value = 13
$ws.Cells.Item($newRow, 6).Value = 100112030
$ws.Cells.Item($newRow, 7).Value = "Poroto granado"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 200
$ws.Cells.Item($newRow, 11).Value = 25000
$ws.Cells.Item($newRow, 12).Value = 27000
$ws.Cells.Item($newRow, 13).Value = 25800
$ws.Cells.Item($newRow, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($newRow, 15).Value = "Región de O'Higgins"
$ws.Cells.Item($newRow, 16).Value = 1032
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
